$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'The applicant possesses several key skills required by the company such as ReactJS, JavaScript, and some level of web development experience. While the applicant''s projects are not directly in line with the company''s requirements, the projects demonstrate a strong grasp of relevant technologies including ReactJS, Django, Python, JavaScript, and Computer Vision. However, the applicant lacks direct experience with MongoDB and NodeJS which are important for this role. Overall, the applicant has a good foundation in relevant technologies and has demonstrated the ability to learn new skills, presenting a strong case for potential success in the role.'

$ws.Range("D2").Value = 4

$ws.Range("E2").Value = 'The applicant demonstrates strong communication skills and a willingness to learn from Japanese work culture, indicating adaptability and a positive attitude towards working in Japan. Moreover, the applicant''s emphasis on seeking feedback, continuous self-improvement, and desire to experience Japanese cultural values align well with the company''s soft skills requirements. The applicant''s clear expression of interest in working in Japan further strengthens the likelihood of successful integration into the work culture.'

$ws.Range("C3").Value = 'The candidate has a strong technical background with experience in web and mobile application development utilizing technologies such as Flutter, Dart, Firebase, ReactJS, NodeJS, HTML, CSS, and WebRTC. Although the applicant''s project experience primarily involves technologies such as Flutter, Dart, and Firebase, the candidate has demonstrated adaptability by working on a variety of projects involving different technology stacks. However, the required skills for the role include MongoDB, ExpressJS, and JavaScript, which are not directly reflected in the candidate''s projects. While the candidate''s projects showcase proficiency in relevant web and mobile technologies, there is a gap in direct experience with the specific technology stack mentioned in the job description. Despite this, the candidate''s diverse skill set and demonstrated ability to learn new technologies provide a strong foundation for success in the role.'

$ws.Range("E3").Value = 'The candidate has demonstrated strong adaptability, communication skills, and a willingness to embrace new challenges, as evidenced by the responses. The candidate''s eagerness to work in Japan, learn a new language, and contribute to cutting-edge AI/ML technologies aligns with the company''s values and culture. Furthermore, the applicant''s proactive approach to addressing weaknesses and commitment to collaborative problem-solving reflect strong teamwork and problem-solving abilities. Overall, the candidate''s responses indicate a high level of enthusiasm, adaptability, and alignment with the company''s soft skill requirements, warranting a score of 4.'

$ws.Range("C4").Value = 'The applicant has strong skills in Python, PyTorch, Django Rest Framework, and ReactJS, which are utilized in the projects ''LLMGuard'', ''Multi Model Data Analysis for Annotation of Human Activities'', and ''Website for the Literature Society of the college''. These projects demonstrate the applicant''s proficiency in web development, machine learning, and backend development. However, the applicant lacks direct experience with MongoDB, ExpressJS, and NodeJS, which are the primary technologies mentioned in the job description. While the applicant''s skills are impressive, there is a gap in the required MERN stack technologies, which may affect the ability to seamlessly transition to the role without additional training or orientation.'

$ws.Range("E4").Value = 'The applicant has exhibited a strong willingness to work in Japan by expressing interest in the AI field and the work culture of Japan. The applicant has shown adaptability by expressing a willingness to work in different countries, including Japan, and has highlighted the desire to learn a new language while working in Japan. Furthermore, the applicant has articulated the expectation of contributing to cutting-edge AI/ML technologies and learning new skills during the internship. The responses indicate good communication skills and a positive attitude towards teamwork and adaptability, aligning with the company''s soft skill requirements.'

$ws.Range("C5").Value = 'The applicant has showcased proficiency in JavaScript, HTML, CSS, and API integration through their projects. These skills align well with the company''s requirements for web development. However, the applicant''s lack of experience in MongoDB, ExpressJS, and NodeJS could be a limiting factor as the job specifically requires knowledge of these technologies. While the applicant''s projects demonstrate strong front-end development skills, the absence of back-end technology experience is a gap that needs to be addressed for this role. It''s crucial for an SDE Intern in MERN stack development to possess comprehensive knowledge across the entire stack, which is currently lacking in the applicant''s skill set.'

$ws.Range("D5").Value = 4

$ws.Range("E5").Value = 'The applicant displayed a strong willingness to work in Japan and adapt to a new culture. Their enthusiasm for Japan''s AI advancements, interest in learning the Japanese language, and admiration for Japanese work culture demonstrate a high level of adaptability and openness. Additionally, the applicant''s emphasis on wanting a stable job, expectations of a decent package, and focus on contributing to cutting-edge AI/ML technologies align with the company''s objectives. The candidate''s proactive approach to overcoming weaknesses and their experience in problem-solving within a team further indicate strong adaptability, teamwork, and problem-solving skills. Based on the provided answers, the applicant exhibits a positive attitude and a genuine interest in working in Japan.'
